$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "FG EMPREIT. MAO DE OBRA LTDA PU_SUL"
$ws.Range("A14").Value = "JAPJ CONSTRUCOES CIVIS LTDA PU_SUDESTE"
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
Write-Output "done"
